# Backlog_3.xlsx edit: the "Semana" column (C) held the text "Semana 03"
# on every data row of both sheets. Replace it with the plain number 3
# (the "Semana 03" shared string becomes unused and is dropped on save).

$wb  = $excel.ActiveWorkbook
$spn = $wb.Worksheets.Item("SPN")
$iti = $wb.Worksheets.Item("ITI")

# ITI!C2 already carries the number style (s=16) that the edited file
# expects everywhere, so borrow its formatting for SPN's Semana column
# before overwriting the values (SPN!C currently uses the text style s=1).
$iti.Range("C2").Copy()
$spn.Range("C2:C27").PasteSpecial(-4122)

$spn.Range("C2:C27").Value = 3
$iti.Range("C2:C22").Value = 3

# Selection / active-sheet bookkeeping to match the saved view state.
$iti.Range("C2:C22").Select()
$spn.Activate()
$spn.Range("C2:C27").Select()
